# Update Overtime (C) and Payment (D) figures on the "Employees" sheet
# for employee rows 2-11, per the fixed-up calculations.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Employees")

$updates = @(
    @{ Row = 2;  C = 122;    D = 2782.5 },
    @{ Row = 3;  C = 128.75; D = 2910 },
    @{ Row = 4;  C = 181.25; D = 4087.5 },
    @{ Row = 5;  C = 119.25; D = 2692.5 },
    @{ Row = 6;  C = 210.5;  D = 4740 },
    @{ Row = 7;  C = 110;    D = 2505 },
    @{ Row = 8;  C = 140.5;  D = 3165 },
    @{ Row = 9;  C = 143;    D = 3217.5 },
    @{ Row = 10; C = 175.75; D = 3960 },
    @{ Row = 11; C = 156.5;  D = 3525 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 3).Value = $u.C
    $ws.Cells.Item($u.Row, 4).Value = $u.D
}
